$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "Check rating"
$ws.Range("C15").Value = "Choose two restaurants on the application and check their rating on Yelp.com`n"
$ws.Range("D15").Value = "The ratings match"
$ws.Range("G15").Value = "Fatih"
$ws.Range("H15").Value = 43172
$ws.Range("I15").Value = "Pass"

$ws.Rows.Item(15).RowHeight = 72

$ws.Range("I16").Select() | Out-Null
